$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

# Rows whose match records were swapped with their neighbour (id/date/B..AC
# payload exchanged while the sequential row-index in column A stays put).
Swap-Rows 68 69
Swap-Rows 81 82
Swap-Rows 83 84
Swap-Rows 90 91
Swap-Rows 102 103
Swap-Rows 108 109

# Row 110 (id 8022181) keeps its own identity but has a handful of odds
# values corrected in place.
$ws.Range("N110").Value2 = 4.75
$ws.Range("P110").Value2 = 1.533
$ws.Range("Q110").Value2 = 1
$ws.Range("R110").Value2 = 1.95
$ws.Range("S110").Value2 = 1.85
